$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELMD")

$ws.Range("H21").Value = 0
$ws.Range("I21").Value = -1200

$ws.Range("D43").Value = 12600
$ws.Range("D44").Value = 4800
$ws.Range("D46").Value = 23300
$ws.Range("D48").Value = 6200
$ws.Range("D49").Value = 1300
$ws.Range("D52").Value = 400
$ws.Range("D54").Value = 27400
$ws.Range("D59").Value = 2900
$ws.Range("D60").Value = 4800
$ws.Range("D66").Value = 4800
$ws.Range("D72").Value = 7600
$ws.Range("D76").Value = 22600

$ws.Range("D91").Value = -500
$ws.Range("E91").Value = -600
$ws.Range("F91").Value = -500
$ws.Range("G91").Value = -500
$ws.Range("I91").Value = -1000
$ws.Range("J91").Value = -800
